# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund holding detail for the quarter)
#    right before the "总计" (total) summary sheet.
# 2) Update the "总计" sheet to add a new first data row summarizing the
#    2022-Q1 figures, shifting the existing rows down by one.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as
# text, even when it looks like a number (needed to keep things like
# fund codes with leading zeros, e.g. "008347", or decimal-looking
# figures such as "4.45" as literal text instead of being silently
# auto-converted to a number by Excel's normal type inference).
# ---------------------------------------------------------------------
function Set-TextValue($sheet, $row, $col, $text) {
    $sheet.Range("ZZ1").Value2 = $text
    $sheet.Range("ZZ1").Copy()
    $sheet.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
}

# =======================================================================
# Step 1: create the new "2022-Q1" worksheet before "总计"
# =======================================================================
# Copying "总计" onto itself inserts a new sheet immediately before it;
# the $totalSheet reference then points at that freshly inserted copy.
$totalSheet.Copy($totalSheet) | Out-Null
$ws = $totalSheet
$ws.Name = "2022-Q1"
$ws.Cells.Clear()

# Re-acquire a handle on the real "总计" sheet (its name is still unique
# since only the copy above was renamed).
$totalSheet = $wb.Worksheets.Item("总计")

# Reuse the existing formatting (bold/centered header style, index style
# on column A) from the "2021-Q4" sheet, which already has the same
# 8-column fund-holding layout.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$template.Range("A2:H4").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$ws.Cells.Item(1, 2).Value2 = "基金代码"
$ws.Cells.Item(1, 3).Value2 = "基金名称"
$ws.Cells.Item(1, 4).Value2 = "基金规模"
$ws.Cells.Item(1, 5).Value2 = "股票总仓位"
$ws.Cells.Item(1, 6).Value2 = "仓位占比"
$ws.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value2 = "仓位排名"

# Prepare the text-forcing helper cell once
$ws.Range("ZZ1").NumberFormat = "@"

# Row 2: 008347 / 中信建投价值甄选混合A
$ws.Cells.Item(2, 1).Value2 = 0
Set-TextValue $ws 2 2 "008347"
Set-TextValue $ws 2 3 "中信建投价值甄选混合A"
Set-TextValue $ws 2 4 "4.45"
Set-TextValue $ws 2 5 "72.24"
Set-TextValue $ws 2 6 "2.51"
Set-TextValue $ws 2 7 "0.1117"
$ws.Cells.Item(2, 8).Value2 = 9

# Row 3: 003822 / 中信建投行业轮换混合A
$ws.Cells.Item(3, 1).Value2 = 1
Set-TextValue $ws 3 2 "003822"
Set-TextValue $ws 3 3 "中信建投行业轮换混合A"
Set-TextValue $ws 3 4 "3.07"
Set-TextValue $ws 3 5 "72.09"
Set-TextValue $ws 3 6 "2.51"
Set-TextValue $ws 3 7 "0.0771"
$ws.Cells.Item(3, 8).Value2 = 8

# Row 4: 003823 / 中信建投行业轮换混合C
$ws.Cells.Item(4, 1).Value2 = 2
Set-TextValue $ws 4 2 "003823"
Set-TextValue $ws 4 3 "中信建投行业轮换混合C"
Set-TextValue $ws 4 4 "0.64"
Set-TextValue $ws 4 5 "72.09"
Set-TextValue $ws 4 6 "2.51"
Set-TextValue $ws 4 7 "0.0161"
$ws.Cells.Item(4, 8).Value2 = 8

$ws.Range("ZZ1").Clear()

# =======================================================================
# Step 2: update the "总计" sheet - add the 2022-Q1 summary row at the
# top of the data (row 2), pushing the existing rows down by one.
# =======================================================================

# Give the new bottom row (row 7) the same index-column style (s=2) as
# the rest of column A before filling it in.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

$summaryRows = @(
    @(0, "2022-Q1", 3, 0.2),
    @(1, "2021-Q4", 5, 0.65),
    @(2, "2021-Q3", 8, 1.08),
    @(3, "2021-Q2", 9, 1.96),
    @(4, "2021-Q1", 5, 0.15),
    @(5, "2020-Q4", 3, 0.04)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = 2 + $i
    $row = $summaryRows[$i]
    $totalSheet.Cells.Item($r, 1).Value2 = $row[0]
    $totalSheet.Cells.Item($r, 2).Value2 = $row[1]
    $totalSheet.Cells.Item($r, 3).Value2 = $row[2]
    $totalSheet.Cells.Item($r, 4).Value2 = $row[3]
}

Write-Host "2022-Q1 sheet added and 总计 updated"
